$d = $word.ActiveDocument

# Locate the exercise-number run ("3. ") at the top of the document so we
# don't have to hard-code absolute character offsets.
$numRng = $d.Content
$found = $numRng.Find.Execute("3. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the '3. ' exercise-number run to update"
}

$digitStart = $numRng.Start
$digitEnd   = $digitStart + 1          # just the "3"
$splitPoint = $digitEnd                # between the digit and ". "

# A temporary bookmark pins the boundary *before* the digit so the run
# holding the leading whitespace doesn't get coalesced with the digit run
# once we rewrite its text.
$d.Bookmarks.Add("zzTempSplit", $d.Range($digitStart, $digitStart))

# Re-seat Word's auto "_GoBack" bookmark right after the new digit; this
# both relocates it here (removing it from wherever it used to live) and
# forces the ". " tail to stay in its own run instead of merging back into
# the word that follows it.
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint))

# 3 -> 4
$d.Range($digitStart, $digitEnd).Text = "4"

# Drop the scaffolding bookmark now that the run split it protected exists.
$d.Bookmarks("zzTempSplit").Delete()
